$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (A5 = 6823) duplicated the StockholdersEquityIncludingPortion /
# AccumulatedOtherComprehensiveIncomeMember combination already present in
# row 2. Remove the duplicate row; everything below shifts up.
$ws.Rows.Item(5).Delete()

# Widen column C to better fit the long dimension-member names
# (126 character-units once saved back to the sheet's <col> width).
$ws.Columns.Item(3).ColumnWidth = 125.28571428571429

$ws.Rows.Item(5).Select()
